$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-10, columns A-F
$data = @(
    @("trainingaudio/02_pitito3.wav", "pngimages/02_pallet.png", "trainingaudio/12_pokika3.wav", "pngimages/12_pie.png", 0.5, -0.5),
    @("trainingaudio/16_kotapi2.wav", "pngimages/16_icecream.png", "trainingaudio/05_titopo2.wav", "pngimages/05_megaphone.png", 0.5, -0.5),
    @("trainingaudio/27_pakapa1.wav", "pngimages/27_kiwi.png", "trainingaudio/26_kapako1.wav", "pngimages/26_pineapple.png", 0.5, -0.5),
    @("trainingaudio/22_kakoki1.wav", "pngimages/22_egg.png", "trainingaudio/14_pokoto1.wav", "pngimages/14_coffee.png", -0.5, 0.5),
    @("trainingaudio/23_patoko1.wav", "pngimages/23_lemon.png", "trainingaudio/12_pokika3.wav", "pngimages/12_pie.png", 0.5, -0.5),
    @("trainingaudio/22_kakoki1.wav", "pngimages/22_egg.png", "trainingaudio/05_titopo2.wav", "pngimages/05_megaphone.png", -0.5, 0.5),
    @("trainingaudio/02_pitito3.wav", "pngimages/02_pallet.png", "trainingaudio/16_kotapi2.wav", "pngimages/16_icecream.png", 0.5, -0.5),
    @("trainingaudio/23_patoko1.wav", "pngimages/23_lemon.png", "trainingaudio/26_kapako1.wav", "pngimages/26_pineapple.png", -0.5, 0.5),
    @("trainingaudio/27_pakapa1.wav", "pngimages/27_kiwi.png", "trainingaudio/14_pokoto1.wav", "pngimages/14_coffee.png", -0.5, 0.5)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
}
